# Shop.xlsx - "fixed for shop pugin"
# Adds a new "Count" (int) column (column I) to the Property1 config sheet,
# mirroring the existing Id/Type/ItemID/... columns: a header label, a type
# row, the Public/Private/Save/Cache/Ref/Upload metadata rows, and a value of
# 1 for every data row (10-71).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: new column header -------------------------------------------------
$ws.Cells.Item(1, 9).Value = "Count"

# --- Rows 2-8: copy column H's formatting onto column I, then fill in the ----
# metadata values (type row + Public/Private/Save/Cache/Ref/Upload flags),
# matching the existing per-row pattern exactly.
$ws.Range("H2:H8").Copy()
$ws.Range("I2:I8").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Cells.Item(2, 9).Value = "int"     # Type row
$ws.Cells.Item(3, 9).Value = $false    # Public
$ws.Cells.Item(4, 9).Value = $false    # Private
$ws.Cells.Item(5, 9).Value = $true     # Save
$ws.Cells.Item(6, 9).Value = $false    # Cache
$ws.Cells.Item(7, 9).Value = $false    # Ref
$ws.Cells.Item(8, 9).Value = $false    # Upload

# --- Rows 10-71: actual data - every item counts as 1 ------------------------
for ($r = 10; $r -le 71; $r++) {
    $ws.Cells.Item($r, 9).Value = 1
}

# Leave the freshly-filled column selected, same as the author did.
$ws.Range("I10:I71").Select()
